$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("precipitation")

# Update existing row 9: Time [s] value changes from 3600 to 500 (Precipitation stays 0)
$ws.Range("A9").Value = 500
$ws.Range("B9").Value = 0

# Add new rows 10-12
$ws.Range("A10").Value = 1000
$ws.Range("B10").Value = 0

$ws.Range("A11").Value = 3000
$ws.Range("B11").Value = 0

$ws.Range("A12").Value = 6000
$ws.Range("B12").Value = 0

# Update the active selection to B2
$ws.Range("B2").Select()
